# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.605.09'
$ws.Range("E2").Value = '  -3.88%  '
$ws.Range("D3").Value = '2.405.97'
$ws.Range("E3").Value = '  -3.72%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = "'511.41"
$ws.Range("E5").Value = '  -4.49%  '
$ws.Range("D6").Value = "'130.12"
$ws.Range("E6").Value = '  -3.05%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").Value = '2.403.19'
$ws.Range("E9").Value = '  -4.44%  '
$ws.Range("D10").Value = "'0.0967"
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = "'4.67"
$ws.Range("E13").Value = '  -10.39%  '
$ws.Range("D14").Value = '2.802.74'
$ws.Range("E14").Value = '  -5.15%  '
$ws.Range("D15").Value = '56.492.65'
$ws.Range("E15").Value = '  -3.97%  '
$ws.Range("D16").Value = "'21.66"
$ws.Range("E16").Value = '  -3.35%  '
$ws.Range("E17").Value = '  -3.23%  '
$ws.Range("D18").Value = '2.378.46'
$ws.Range("E18").Value = '  -5.43%  '
$ws.Range("D19").Value = "'10.26"
$ws.Range("E19").Value = '  -3.54%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'313.46"
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = "'4.07"
$ws.Range("E21").Value = '  -4.42%  '
$ws.Range("D22").Value = "'6.32"
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = "'65.64"
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '2.486.27'
$ws.Range("E26").Value = '  -5.53%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = "'0.152"
$ws.Range("E27").Value = '  -4.56%  '
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").Value = "'0.376"
$ws.Range("E28").Value = '  -8.46%  '
$ws.Range("D29").Value = "'7.25"
$ws.Range("D30").Value = "'174.98"
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").Value = '  -5.56%  '
$ws.Range("D34").Value = "'1.11"
$ws.Range("E34").Value = '  -6.67%  '
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = "'0.994"
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -2.39%  '
$ws.Range("D38").Value = "'1.21"
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = '  -5.77%  '
$ws.Range("D40").Value = "'35.87"
$ws.Range("E40").Value = '  -1.77%  '
$ws.Range("D41").Value = "'1.44"
$ws.Range("E41").Value = '  -5.18%  '
$ws.Range("D42").Value = "'0.789"
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("D43").Value = "'134.65"
$ws.Range("E43").Value = '  +1.95%  '
$ws.Range("E44").Value = '  -4.04%  '
$ws.Range("D45").Value = "'4.91"
$ws.Range("E45").Value = '  -3.06%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").Value = "'256.48"
$ws.Range("E46").Value = '  -7.50%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'0.573"
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("D48").Value = "'0.0901"
$ws.Range("E48").Value = '  -3.58%  '
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("E50").Value = '  -5.00%  '
$ws.Range("D51").Value = "'16.79"
$ws.Range("E51").Value = '  -5.11%  '
